# Refactor the "synthetic array" of status icon/label shared strings:
#   black square (⬛ / noir)  -> blue book  (📘 / bleu)
#   red   square (🟥)         -> red  book  (📕)
#   orange square (🟧)        -> orange book (📙)
#   green square (🟩)        -> green book (📗)
# "rouge", "orange", "vert" text labels are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "⬛"   = "📘"
    "🟥"   = "📕"
    "🟧"   = "📙"
    "🟩"   = "📗"
    "noir" = "bleu"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range($col + $r)
        $val = $cell.Value()
        if ($null -ne $val -and $map.ContainsKey([string]$val)) {
            $cell.Value = $map[[string]$val]
        }
    }
}
